$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2600
$ws.Range("I64").Value = 2514.2856
$ws.Range("K64").Value = 2514.2856
$ws.Range("M64").Value = -2266.2856
$ws.Range("H67").Value = 2600
$ws.Range("I67").Value = 2514.2856
$ws.Range("K67").Value = 2514.2856
$ws.Range("M67").Value = -1656.2856
$ws.Range("H69").Value = 1745
$ws.Range("I69").Value = 1490
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 4470
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = -3596
$ws.Range("N69").Value = -7748
$ws.Range("H70").Value = 1528.8
$ws.Range("I70").Value = 1715.3846
$ws.Range("J70").Value = 1326.6666
$ws.Range("K70").Value = 5146.1538
$ws.Range("L70").Value = 3979.9998
$ws.Range("M70").Value = -4876.1538
$ws.Range("N70").Value = -4519.9998
$ws.Range("H72").Value = 1745
$ws.Range("I72").Value = 1490
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 13410
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = -9042
$ws.Range("N72").Value = -26736
$ws.Range("H73").Value = 1528.8
$ws.Range("I73").Value = 1715.3846
$ws.Range("J73").Value = 1326.6666
$ws.Range("K73").Value = 5146.1538
$ws.Range("L73").Value = 3979.9998
$ws.Range("M73").Value = -4210.1538
$ws.Range("N73").Value = -5851.9998
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 900
$ws.Range("I41").Value = 900
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -486
$ws.Range("H61").Value = 11112826
$ws.Range("I61").Value = 13890632
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 13890632
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -13890420
$ws.Range("N61").Value = -2024
$ws.Range("H132").Value = 2675610
$ws.Range("I132").Value = 1458.6
$ws.Range("J132").Value = 8405934
$ws.Range("K132").Value = 4375.799999999999
$ws.Range("L132").Value = 25217802
$ws.Range("M132").Value = -1845.799999999999
$ws.Range("N132").Value = -25222862
$ws.Range("H136").Value = 11112826
$ws.Range("I136").Value = 13890632
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 41671896
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -41669346
$ws.Range("N136").Value = -9900
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1780.909
$ws.Range("I99").Value = 1390.4166
$ws.Range("J99").Value = 2822.2222
$ws.Range("K99").Value = 1390.4166
$ws.Range("L99").Value = 2822.2222
$ws.Range("M99").Value = 107.5834
$ws.Range("N99").Value = -5818.2222
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1100
$ws.Range("I6").Value = 200
$ws.Range("K6").Value = 200
$ws.Range("M6").Value = -87
$ws.Range("H33").Value = 13219
$ws.Range("I33").Value = 1438
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 1438
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -1059
$ws.Range("N33").Value = -25758
$ws.Range("H39").Value = 3250
$ws.Range("I39").Value = 3250
$ws.Range("K39").Value = 3250
$ws.Range("M39").Value = -2859
$ws.Range("H49").Value = 3250
$ws.Range("I49").Value = 3250
$ws.Range("K49").Value = 3250
$ws.Range("M49").Value = -3068
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 27783988
$ws.Range("I5").Value = 35714790
$ws.Range("J5").Value = 26175
$ws.Range("K5").Value = 107144370
$ws.Range("L5").Value = 78525
$ws.Range("M5").Value = -107144258
$ws.Range("N5").Value = -78749
$ws.Range("H107").Value = 40004690
$ws.Range("I107").Value = 170.35294
$ws.Range("J107").Value = 125014290
$ws.Range("K107").Value = 511.05882
$ws.Range("L107").Value = 375042870
$ws.Range("M107").Value = 1408.94118
$ws.Range("N107").Value = -375046710
$ws.Range("H131").Value = 888.62
$ws.Range("J131").Value = 897.5714
$ws.Range("L131").Value = 2692.7142
$ws.Range("N131").Value = -12772.7142
$ws.Range("H132").Value = 20005208
$ws.Range("I132").Value = 835.3333
$ws.Range("J132").Value = 31257668
$ws.Range("K132").Value = 7517.9997
$ws.Range("L132").Value = 281319012
$ws.Range("M132").Value = -4987.9997
$ws.Range("N132").Value = -281324072
$ws.Range("H133").Value = 83335770
$ws.Range("I133").Value = 83335770
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 250007310
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -250002250
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 15154281
$ws.Range("I134").Value = 29412430
$ws.Range("K134").Value = 88237290
$ws.Range("M134").Value = -88232220
$ws.Range("H135").Value = 27783988
$ws.Range("I135").Value = 35714790
$ws.Range("J135").Value = 26175
$ws.Range("K135").Value = 321433110
$ws.Range("L135").Value = 235575
$ws.Range("M135").Value = -321430575
$ws.Range("N135").Value = -240645
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H61").Value = 1415.9412
$ws.Range("J61").Value = 1307.6364
$ws.Range("L61").Value = 1307.6364
$ws.Range("N61").Value = -1711.6364
$ws.Range("H68").Value = 1539.0769
$ws.Range("I68").Value = 1500.6666
$ws.Range("J68").Value = 1625.5
$ws.Range("K68").Value = 1500.6666
$ws.Range("L68").Value = 1625.5
$ws.Range("M68").Value = -751.6666
$ws.Range("N68").Value = -3123.5
$ws.Range("H71").Value = 1539.0769
$ws.Range("I71").Value = 1500.6666
$ws.Range("J71").Value = 1625.5
$ws.Range("K71").Value = 7503.333000000001
$ws.Range("L71").Value = 8127.5
$ws.Range("M71").Value = -3759.333000000001
$ws.Range("N71").Value = -15615.5
$ws.Range("H113").Value = 1415.9412
$ws.Range("J113").Value = 1307.6364
$ws.Range("L113").Value = 1307.6364
$ws.Range("N113").Value = -5647.6364
$ws.Range("H132").Value = 8147.7856
$ws.Range("I132").Value = 2624.111
$ws.Range("J132").Value = 12290.542
$ws.Range("K132").Value = 7872.333
$ws.Range("L132").Value = 36871.626
$ws.Range("M132").Value = -5342.333
$ws.Range("N132").Value = -41931.626
